$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.9168783282299557,0.1103637748559478,0.1212360588839161,0.1322696356982007,1.940048721761215,0,0.07973214163530429,0,0.1762184897145787,0.478284596782089,0.2631712867810023,0,2.252590438682958,5.132683292418534),
    @(0.8797396025435944,0.1093442579306725,0.119489934097345,0.1323217267412389,1.946675811283775,0,0.07973214163530429,0,0.1771558591426476,0.4443327767251333,0.2593479686275444,0,2.272981447585809,5.157466130098101),
    @(0.8572861856528391,0.1087140804762328,0.1184611853797151,0.13239750608839,1.951551618945011,0,0.07973214163530429,0,0.1777840706764486,0.4236228971695937,0.2571021721292723,0,2.286138749727508,5.17483604580633),
    @(0.8482249817101319,0.1084562410397041,0.1180529336918283,0.1324394302075298,1.953741668300538,0,0.07973214163530429,0,0.1780533363988273,0.4152183855656517,0.2562126944636702,0,2.291660731319653,5.182456355751626),
    @(0.8467257580251157,0.1084133648533268,0.1179858083541632,0.1324470595754459,1.954117599665921,0,0.07973214163530429,0,0.1780988495502118,0.4138249492640256,0.2560665532878588,0,2.292587334124478,5.183754448296312),
    @(0.8571636227548538,0.1087106073412869,0.1184556350456418,0.1323980267362312,1.951580331961061,0,0.07973214163530429,0,0.1777876483517176,0.4235094087986226,0.2570900720931206,0,2.286212572078821,5.174936621062415),
    @(0.904000713041512,0.1100131246874128,0.1206250292396334,0.1322785252058924,1.942166452186058,0,0.07973214163530429,0,0.1765307776214335,0.4665499697962332,0.2618319617224785,0,2.259489058073968,5.140781823717219),
    @(0.9985964627635724,0.1125334978329988,0.1252209016906392,0.1323904215224907,1.930096956476262,0,0.07973214163530429,0,0.1744829630841593,0.5520164383819974,0.271933434149048,0,2.212136269059181,5.090868847179252),
    @(1.06974152430513,0.1143639177496993,0.1288028934697252,0.1326821564595591,1.925113754239973,0,0.07973214163530429,0,0.1732313271048405,0.6154363994872654,0.2798391457279337,0,2.180420048271301,5.064577286381649),
    @(1.102458693753107,0.1151918697636631,0.1304764475679718,0.1328600511805718,1.923687794362806,0,0.07973214163530429,0,0.1727165747196011,0.6444200335887729,0.2835397588262794,0,2.166657716816241,5.054865695375867),
    @(1.114897916969738,0.1155047009853192,0.1311164599157166,0.1329338835726332,1.923268521099914,0,0.07973214163530429,0,0.1725294849773746,0.6554141001142,0.2849559697561546,0,2.161541927891806,5.051511087788811),
    @(1.112216699417075,0.1154373584393582,0.1309783438316288,0.1329176952743474,1.923353454085117,0,0.07973214163530429,0,0.1725694299009675,0.6530455140275251,0.2846503040648685,0,2.162639447867365,5.052219204819892),
    @(1.103481079871159,0.1152176206345885,0.1305289764119522,0.1328659959672436,1.923650882597514,0,0.07973214163530429,0,0.1727010258002331,0.6453241529427203,0.283655974064601,0,2.166234919632505,5.054583239344538),
    @(1.098136740470835,0.1150829336747563,0.1302545409574094,0.1328351700292636,1.923848779300059,0,0.07973214163530429,0,0.1727826519830771,0.6405969994779639,0.2830488509347759,0,2.168449713990833,5.056073326548159),
    @(1.067610412936375,0.1143097126997148,0.1286944042891918,0.1326714372757856,1.925223847051058,0,0.07973214163530429,0,0.1732660649201989,0.6135448861811312,0.2795993894300608,0,2.181332841567282,5.065257176133031),
    @(1.048973316645146,0.1138341460382009,0.1277485559645157,0.1325825438171222,1.926282632845414,0,0.07973214163530429,0,0.1735765998907759,0.5969830655719761,0.2775098730280803,0,2.189406684822844,5.071466812294716),
    @(1.038287019046777,0.1135601699437743,0.1272086827536896,0.1325356692679698,1.926970783114811,0,0.07973214163530429,0,0.1737603539015247,0.5874697539454985,0.2763178598222282,0,2.194113188930384,5.075250099670086),
    @(1.034674560621227,0.1134673308606651,0.1270266063821595,0.1325205299924548,1.927217382201761,0,0.07973214163530429,0,0.1738234536614875,0.584250894515975,0.2759159552434909,0,2.195717489279043,5.076567421735035),
    @(1.050953830183204,0.1138848168735933,0.1278488137177476,0.132591566556421,1.926161731614314,0,0.07973214163530429,0,0.1735430108200298,0.5987447990603698,0.2777312900413307,0,2.188540727086451,5.070783881746138),
    @(1.106045594144092,0.1152821819946013,0.130660796775274,0.1328810060171293,1.923560246414453,0,0.07973214163530429,0,0.172662160383446,0.6475916041456742,0.2839476303361721,0,2.16517624490422,5.053880103055377),
    @(1.142341966614083,0.1161913742745426,0.1325351277155562,0.1331078553315663,1.922563496479924,0,0.07973214163530429,0,0.1721321381179273,0.679623905565137,0.2880969787200769,0,2.150464064510023,5.044714795463619),
    @(1.122943575962438,0.1157064993367811,0.1315314406662083,0.1329833433204222,1.923031184679459,0,0.07973214163530429,0,0.1724108489036347,0.6625179763570372,0.2858745095444846,0,2.158265185841087,5.049434388540476),
    @(1.050058351434501,0.11386191035335,0.1278034750122004,0.1325874741933752,1.926216143593393,0,0.07973214163530429,0,0.1735581801672552,0.5979482934514806,0.2776311585020892,0,2.188932024902764,5.071091970246812),
    @(0.9727146670326192,0.1118553630851693,0.1239413156595432,0.1323232439733708,1.932679098916573,0,0.07973214163530429,0,0.1749924496055399,0.5287837424394013,0.271933434149048,0,2.224406232196659,5.102547192352517)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, 2 + $j).Value = $row[$j]
    }
}
